$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.397.49"
$ws.Range("E2").Value = "'  -1.42%  "
$ws.Range("D3").Value = "'2.981.79"
$ws.Range("E3").Value = "'  -0.42%  "
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("D5").Value = "'505.60"
$ws.Range("E5").Value = "'  +0.72%  "
$ws.Range("D6").Value = "'137.21"
$ws.Range("E6").Value = "'  -1.11%  "
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E8").Value = "'  -0.75%  "
$ws.Range("D9").Value = "'7.19"
$ws.Range("E9").Value = "'  -1.67%  "
$ws.Range("E10").Value = "'  -0.99%  "
$ws.Range("E11").Value = "'  +1.05%  "
$ws.Range("D12").Value = "'3.490.02"
$ws.Range("E12").Value = "'  -0.50%  "
$ws.Range("D13").Value = "'0.126"
$ws.Range("E13").Value = "'  -1.06%  "
$ws.Range("D14").Value = "'25.80"
$ws.Range("E14").Value = "'  -1.82%  "
$ws.Range("D15").Value = "'0.0000163"
$ws.Range("E15").Value = "'  +1.27%  "
$ws.Range("D16").Value = "'56.314.46"
$ws.Range("E16").Value = "'  -1.66%  "
$ws.Range("D17").Value = "'2.978.55"
$ws.Range("E17").Value = "'  -0.63%  "
$ws.Range("D18").Value = "'5.98"
$ws.Range("E18").Value = "'  -1.92%  "
$ws.Range("D19").Value = "'12.86"
$ws.Range("E19").Value = "'  +1.29%  "
$ws.Range("D20").Value = "'8.06"
$ws.Range("E20").Value = "'  +2.01%  "
$ws.Range("D21").Value = "'331.51"
$ws.Range("E21").Value = "'  +3.13%  "
$ws.Range("E22").Value = "'  +0.23%  "
$ws.Range("E23").Value = "'  -0.40%  "
$ws.Range("D24").Value = "'64.52"
$ws.Range("E24").Value = "'  +1.50%  "
$ws.Range("D25").Value = "'3.106.81"
$ws.Range("E25").Value = "'  -0.41%  "
$ws.Range("E26").Value = "'  -0.09%  "
$ws.Range("E27").Value = "'  -0.49%  "
$ws.Range("D28").Value = "'0.0₃0917"
$ws.Range("E28").Value = "'  +1.91%  "
$ws.Range("D29").Value = "'6.35"
$ws.Range("E29").Value = "'  -3.30%  "
$ws.Range("D30").Value = "'6.93"
$ws.Range("E30").Value = "'  -3.26%  "
$ws.Range("E31").Value = "'  +0.00%  "
$ws.Range("B32").Value = "'EthereumClassic"
$ws.Range("C32").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'20.20"
$ws.Range("E32").Value = "'  -0.60%  "
$ws.Range("B33").Value = "'Fetch.AI"
$ws.Range("C33").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.15"
$ws.Range("E33").Value = "'  -1.71%  "
$ws.Range("D34").Value = "'152.91"
$ws.Range("E34").Value = "'  -1.43%  "
$ws.Range("D35").Value = "'4.46"
$ws.Range("E35").Value = "'  -2.68%  "
$ws.Range("D36").Value = "'5.78"
$ws.Range("E36").Value = "'  -0.28%  "
$ws.Range("D37").Value = "'26.12"
$ws.Range("E37").Value = "'  +7.03%  "
$ws.Range("E38").Value = "'  -0.80%  "
$ws.Range("D39").Value = "'0.0659"
$ws.Range("E39").Value = "'  -1.08%  "
$ws.Range("D40").Value = "'3.017.60"
$ws.Range("E40").Value = "'  -0.29%  "
$ws.Range("D41").Value = "'36.90"
$ws.Range("E41").Value = "'  -2.41%  "
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "'  -0.23%  "
$ws.Range("D43").Value = "'3.79"
$ws.Range("E43").Value = "'  +0.61%  "
$ws.Range("E44").Value = "'  +0.71%  "
$ws.Range("D45").Value = "'2.178.98"
$ws.Range("E45").Value = "'  -1.27%  "
$ws.Range("E46").Value = "'  -3.24%  "
$ws.Range("E47").Value = "'  -2.84%  "
$ws.Range("E48").Value = "'  -2.75%  "
$ws.Range("D49").Value = "'0.0234"
$ws.Range("E49").Value = "'  -0.68%  "
$ws.Range("D50").Value = "'19.49"
$ws.Range("E50").Value = "'  +0.66%  "
$ws.Range("E51").Value = "'  -3.00%  "
